# Insert a new weekly price-record row for "Arveja Verde" (Feria Lagunitas de
# Puerto Montt) right after the existing row 116, pushing every following
# record down by one row (old row 117 becomes 118, ..., old row 144 becomes
# row 145). This mirrors the target diff, which grows the used range from
# A1:R144 to A1:R145 and shifts all data rows 117-144 down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 117; Excel shifts rows 117..144 down to
# 118..145 and the new row inherits formatting (incl. the date style) from
# the row above it, matching style s="2" seen on column D in the diff.
$ws.Rows.Item(117).Insert()

# Populate the newly inserted row 117 with the new record's data.
$ws.Range("A117").Value2 = 4
$ws.Range("B117").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C117").Value2 = "Los Lagos"
$ws.Range("D117").Value2 = 44855
$ws.Range("E117").Value2 = 10
$ws.Range("F117").Value2 = 100112022
$ws.Range("G117").Value2 = "Arveja Verde"
$ws.Range("H117").Value2 = "Perfection"
$ws.Range("I117").Value2 = "Primera"
$ws.Range("J117").Value2 = 70
$ws.Range("K117").Value2 = 27000
$ws.Range("L117").Value2 = 27000
$ws.Range("M117").Value2 = 27000
$ws.Range("N117").Value2 = "`$/malla 25 kilos"
$ws.Range("O117").Value2 = "Provincia de Huasco"
$ws.Range("P117").Value2 = 1080
$ws.Range("Q117").Value2 = 25
$ws.Range("R117").Value2 = "Hortaliza"
